# feat: add 2022-Q4 data
#
# 1) Insert a brand-new worksheet "2022-Q4" right before the existing
#    "2022-Q3" sheet and fill it with the new quarter's fund-holding table.
# 2) Insert a new top data row in the "总计" (summary) sheet for the
#    2022-Q4 totals, pushing the existing quarter rows down by one (and
#    renumbering the row-index column to match).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (row 1), columns B..H - matches the other quarter sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Range($headerCols[$i] + "1")
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data rows 2..6. Columns B-G are stored as literal text (to preserve
# formatting such as leading zeros in fund codes and trailing zeros in
# percentages), column H (rank) is numeric, same as the other sheets.
$rows = @(
    @("161724", "招商中证煤炭等权指数（LOF）A", "17.24", "93.84", "3.08", "0.5310", 10),
    @("009837", "华夏磐锐一年定期开放混合A",     "14.15", "75.21", "3.49", "0.4938", 7),
    @("013596", "招商中证煤炭等权指数（LOF）C", "1.56",  "93.84", "3.08", "0.0480", 10),
    @("009838", "华夏磐锐一年定期开放混合C",     "0.39",  "75.21", "3.49", "0.0136", 7),
    @("016347", "招商中证煤炭等权指数（LOF）E", "0.20",  "93.84", "3.08", "0.0062", 10)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $aCell = $q4.Range("A" + $rowNum)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $bCell = $q4.Range("B" + $rowNum)
    $bCell.NumberFormat = "@"
    $bCell.Value = $data[0]

    $cCell = $q4.Range("C" + $rowNum)
    $cCell.NumberFormat = "@"
    $cCell.Value = $data[1]

    $dCell = $q4.Range("D" + $rowNum)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[2]

    $eCell = $q4.Range("E" + $rowNum)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[3]

    $fCell = $q4.Range("F" + $rowNum)
    $fCell.NumberFormat = "@"
    $fCell.Value = $data[4]

    $gCell = $q4.Range("G" + $rowNum)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data[5]

    $q4.Range("H" + $rowNum).Value = $data[6]
}

# ---------------------------------------------------------------------
# Step 2: new top row in "总计" for the 2022-Q4 totals
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Newly inserted row inherits formatting from the row above (the header);
# clear it on B2:D2 so it matches the plain data rows below.
$summary.Range("B2:D2").ClearFormats()

$a2 = $summary.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 1.09

# The row-index column (A) is a plain 0-based sequence; renumber the rows
# that got pushed down so it stays contiguous (0,1,2,3,4).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
